# Append two new data rows (24 and 25) to Sheet1, mirroring the existing
# daily station-pair pattern: one row for 四方坪站 (station 1) and one for
# 高岭站 (station 2) on date serial 46003 (2025-12-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24 - 四方坪站
$ws.Cells.Item(24, 1).Value = 46003
$ws.Cells.Item(24, 2).Value = "四方坪站"
$ws.Cells.Item(24, 3).Value = 8938.24
$ws.Cells.Item(24, 4).Value = 7604.64
$ws.Cells.Item(24, 5).Value = 2951.32
$ws.Cells.Item(24, 6).Value = 412

# Row 25 - 高岭站
$ws.Cells.Item(25, 1).Value = 46003
$ws.Cells.Item(25, 2).Value = "高岭站"
$ws.Cells.Item(25, 3).Value = 4923.18
$ws.Cells.Item(25, 4).Value = 4286.1400000000003
$ws.Cells.Item(25, 5).Value = 1228.83
$ws.Cells.Item(25, 6).Value = 180

# Match the updated viewport: scroll so row 19 is at the top and select G25.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G25").Select()
